# Updates crypto price/volume data per the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.347.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.244.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.18%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.85%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.233.45"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.583"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -8.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "632.46"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.90%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.51"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.761.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.157.06"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.73%  "

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.71"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.61%  "

# Row 19
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.116"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.258.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.85%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -8.40%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.39%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.51"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.18%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -7.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.45"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.66"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.22"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.98"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.88%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.25"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.81%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "560.54"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +12.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.69%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.28"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.59%  "

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.01%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.590.22"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.36%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "CoreDAO"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.71"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +16.91%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.01%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.34%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0703"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.127"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.08%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "31.81"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.20%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.333"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.30%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.26"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0410"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.51%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.56"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.66%  "
